$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 158, shifting existing rows 158-167 down to 159-168.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row with the new weekly price-report entry.
$ws.Range("A158").Value = 11
$ws.Range("B158").Value = "Vega Monumental Concepción"
$ws.Range("C158").Value = "Bíobío"
$ws.Range("D158").Value = 44516
$ws.Range("E158").Value = 8
$ws.Range("F158").Value = 100112017
$ws.Range("G158").Value = "Apio"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 550
$ws.Range("K158").Value = 5500
$ws.Range("L158").Value = 6000
$ws.Range("M158").Value = 5773
$ws.Range("N158").Value = "$/docena de matas"
$ws.Range("O158").Value = "Región de Coquimbo"
$ws.Range("P158").Value = 962
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"
